# Append two new data rows (for the tensorflow/ranking repository) to the
# "dataset" worksheet, growing it from 172 to 174 data rows.
#
# The new rows are copied from the last existing row so they inherit its
# formatting/styles, then only the cells that actually differ from that
# source row (the id in column A, and OS in column J for the first of the
# two new rows) are corrected.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dataset")

$lastRow = 172
$newRow1 = $lastRow + 1
$newRow2 = $lastRow + 2

$srcRange = "A" + $lastRow + ":O" + $lastRow

# Row 173 (id 172): same as the source row except OS (column J) is 0.
$ws.Range($srcRange).Copy($ws.Range("A" + $newRow1 + ":O" + $newRow1))
$ws.Cells.Item($newRow1, 1).Value = 172
$ws.Cells.Item($newRow1, 10).Value = "'0"
$ws.Cells.Item($newRow1, 10).Style = $ws.Cells.Item($newRow1, 9).Style

# Row 174 (id 173): identical flags to the source row, only the id changes.
$ws.Range($srcRange).Copy($ws.Range("A" + $newRow2 + ":O" + $newRow2))
$ws.Cells.Item($newRow2, 1).Value = 173

$excel.CutCopyMode = $false
